$d = $word.ActiveDocument

# 1) Update the date in the first line (simple find/replace)
$null = $d.Content.Find.Execute("30.08.24", $true, $false, $false, $false, $false, $true, 1, $false, "22.08.24", 2)

# 2) Replace the paper title (2nd paragraph)
$null = $d.Content.Find.Execute("Platypus: A Generalized Specialist Model for Reading Text in Various Forms", $true, $false, $false, $false, $false, $true, 1, $false, "Approaching Deep Learning through the Spectral Dynamics of Weights", 2)

# 3) Replace the three body paragraphs with their new content (whole-paragraph rewrite)
$d.Paragraphs.Item(3).Range.Text = "היום נסקור מאמר החוקר מה הסיבות לתופעה של גרוקינג. למי שלא מכיר גרוקינג זו תופעה די מעניינת המתרחשת כאשר ממשיכים לאמן רשת נוירונים (למרות שזה קורה גם במודלים אחרים) גם אחרי לוס הוולידציה מתחיל לעלות (כלומר אנו נכנסים למשטר אוורפיט). מתברר אם לא עוצרים וממשיכים לאמן לוס הוולידציה מתחיל לרדת כלומר המודל נכנס למשטר ההכללה כלומר לומד את ה״חוקיות האמיתית״ מאחורי הדאטה. "
$d.Paragraphs.Item(4).Range.Text = "התופעה הזו היא מקרה פרטי של double descent (יש גם multiple descent) שמתרחש גם אם אנו מוסיפים פרמטרים למודל בצורה עקבית ומגיעים למצב שיש לנו over-parametrization. כלומר יש המודל שלנו לכאורה מתחיל ״יותר מדי פרמטרים״ כדי ״להבין את הדאטה״. וגם שם זה קורה בצורה בלתי רציפה כלומר יש אינטרוול של פרמטרים שביצועי המודל יורדים עבורם ורק אז מתחילים לרדת. "
$d.Paragraphs.Item(5).Range.Text = "המאמר חוקר מה קורה עם משקלי המודל כאשר הוא נכנס למשטר הגרוקינג. מתברר שתופעה הגרוקינג קשורה לירידה בראנק של מטריצות המשקלים של המודל. בשבילי זה די אינטואיטיבי כי לדעתי במהלך גרוקינג המודל מצליח להתכנס ל״פתרון פשוט ביותר עבור הדאטהסט. פתרון פשוט הכוונה הוא מודל שאפקטיבית הוא קטן, כלומר רוב וקטורי המשקלים בו או אפס או תלוים לינארית זה בזה."

# 4) Delete the closing "זה כל הפרטים..." paragraph entirely (paragraph 6)
$d.Paragraphs.Item(6).Range.Delete()

# 5) Update the arxiv link (now the last paragraph)
$null = $d.Content.Find.Execute("2408.14805", $true, $false, $false, $false, $false, $true, 1, $false, "2408.11804", 2)

